# Add 2022-Q3 data
#
# 1. Insert a new row at the top of the "总计" (summary) sheet's data
#    block and fill it with the 2022-Q3 totals; renumber the existing
#    index column.
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计"
#    (so the tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, ...)
#    and populate it with the per-fund holdings for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 总计 sheet: insert the 2022-Q3 summary row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$a2 = $summary.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.09

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet with per-fund holdings, placed right after
#    "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3Rows = @(
    @(0, "006969", "圆信永丰高端制造混合",   "0.91", "87.79", "4.26", "0.0388", 2),
    @(1, "009847", "圆信永丰研究精选混合A",  "1.14", "89.51", "3.30", "0.0376", 9),
    @(2, "009848", "圆信永丰研究精选混合C",  "0.44", "89.51", "3.30", "0.0145", 9),
    @(3, "009054", "圆信永丰沣泰混合",       "0.23", "26.81", "1.28", "0.0029", 5)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").NumberFormat = "@"
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").NumberFormat = "@"
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").NumberFormat = "@"
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").NumberFormat = "@"
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").NumberFormat = "@"
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r++
}
